$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of B2:H2 (row 2's customer/order data), leaving the
# cells empty while preserving their existing formatting/style.
$ws.Range("B2:H2").ClearContents()
